$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The weekly update inserts a new price observation as row 101 (pushing the
# existing rows 101-108 down to 102-109), matching the "Fruta / hortaliza,
# semanal" refresh pattern used across these sheets.
$ws.Rows.Item(101).Insert()

$ws.Range("A101").Value = 8
$ws.Range("B101").Value = "Terminal La Palmera de La Serena"
$ws.Range("C101").Value = "Coquimbo"
$ws.Range("D101").Value = 44783
$ws.Range("E101").Value = 4
$ws.Range("F101").Value = 100112052
$ws.Range("G101").Value = "Albahaca"
$ws.Range("H101").Value = "Sin especificar"
$ws.Range("I101").Value = "Primera"
$ws.Range("J101").Value = 1600
$ws.Range("K101").Value = 3300
$ws.Range("L101").Value = 3500
$ws.Range("M101").Value = 3400
$ws.Range("N101").Value = "$/paquete"
$ws.Range("O101").Value = "Región de Arica y Parinacota"
$ws.Range("P101").Value = 3400
$ws.Range("Q101").Value = 1
$ws.Range("R101").Value = "Hortaliza"
